# Add 2022-Q3 data
#
# 1) Insert a new row into the "总计" (totals) sheet for the 2022-Q3 quarter,
#    pushing the existing 2022-Q1 / 2021-Q3 rows down.
# 2) Insert a new worksheet "2022-Q3" (with fund-level detail) right after
#    "总计" and before "2022-Q1".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: "总计" sheet - insert new summary row for 2022-Q3
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("B2:D2").ClearFormats()

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.09

# Re-sequence the index column (A) for the rows that shifted down
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2

# Match column A's style (bold/centered/bordered) used by the other rows
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Step 2: create the "2022-Q3" detail sheet, positioned right after "总计"
# ---------------------------------------------------------------------------
$wsQ1 = $wb.Worksheets.Item("2022-Q1")

# Duplicate the "2022-Q1" sheet (so formatting/styles match the other
# quarterly detail sheets) and place the copy right after "总计".
$wsQ1.Copy($null, $wsTotal)
$wsQ3 = $wb.Worksheets.Item("2022-Q1 (2)")
$wsQ3.Name = "2022-Q3"

# The source sheet has 3 data rows; 2022-Q3 only needs 2, so drop the last one.
$wsQ3.Rows.Item(4).Delete()

$wsQ3.Range("A2").Value = 0
$wsQ3.Range("B2").Value = "'011160"
$wsQ3.Range("C2").Value = "富国质量成长6个月持有期混合A"
$wsQ3.Range("D2").Value = "'3.70"
$wsQ3.Range("E2").Value = "'85.89"
$wsQ3.Range("F2").Value = "'2.43"
$wsQ3.Range("G2").Value = "'0.0899"
$wsQ3.Range("H2").Value = 10

$wsQ3.Range("A3").Value = 1
$wsQ3.Range("B3").Value = "'011161"
$wsQ3.Range("C3").Value = "富国质量成长6个月持有期混合C"
$wsQ3.Range("D3").Value = "'0.14"
$wsQ3.Range("E3").Value = "'85.89"
$wsQ3.Range("F3").Value = "'2.43"
$wsQ3.Range("G3").Value = "'0.0034"
$wsQ3.Range("H3").Value = 10

# Setting text-like values ("011160", "3.70", ...) makes Excel apply an
# explicit "@" text number-format; strip that so the cells match the plain
# (unformatted) text cells used elsewhere in the workbook.
$wsQ3.Range("B2:G3").ClearFormats()

# ---------------------------------------------------------------------------
# Restore the originally active sheet/tab ("2021-Q3"), since adding/copying
# sheets shifts Excel's active-sheet selection.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q3").Activate()
